$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Centroid Calculation block (columns G:I, rows 29-34) ---
# (string values are entered in the same order the author created them so the
#  shared-string table indices line up with the target workbook)

# Row 29 - title
$ws.Range("G29").Value = "Centroid Calculation"

# Row 30 - INDEX row
$ws.Range("G30").Value = "INDEX"
$ws.Range("H30").Value = 2
$ws.Range("I30").Value = 5000

# Row 31 - i_l row
$ws.Range("G31").Value = "i_l"
$ws.Range("H31").Formula = "=H30-1"
$ws.Range("I31").Value = 1

# Row 32 - i_r row
$ws.Range("G32").Value = "i_r"
$ws.Range("H32").Formula = "=H30+1"
$ws.Range("I32").Value = 5000

# Row 29 - Value column header
$ws.Range("I29").Value = "Value"

# Row 34 - centroid result
$ws.Range("G34").Value = "centroid"
$ws.Range("H34").Formula = "=((H30*I30)+(H31*I31)+(H32*I32))/(SUM(I30:I32))"

# --- Update selection to match author's final cursor position ---
$ws.Range("I31").Select()

# Best-effort: scroll the window so row 21 is at the top (topLeftCell="A21").
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
